$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 2990
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H17").Value = 2024.7059
$ws.Range("J17").Value = 2024.7059
$ws.Range("L17").Value = 6074.1177
$ws.Range("N17").Value = -6410.1177
$ws.Range("H19").Value = 17628954
$ws.Range("I19").Value = 15651806
$ws.Range("J19").Value = 20001532
$ws.Range("K19").Value = 15651806
$ws.Range("L19").Value = 20001532
$ws.Range("M19").Value = -15651631
$ws.Range("N19").Value = -20001882
$ws.Range("H20").Value = 2783
$ws.Range("I20").Value = 2783
$ws.Range("K20").Value = 2783
$ws.Range("M20").Value = -2553
$ws.Range("H33").Value = 87.111115
$ws.Range("I33").Value = 92
$ws.Range("K33").Value = 92
$ws.Range("M33").Value = 137
$ws.Range("H35").Value = 2783
$ws.Range("I35").Value = 2783
$ws.Range("K35").Value = 2783
$ws.Range("M35").Value = -2404
$ws.Range("H39").Value = 550.2857
$ws.Range("I39").Value = 85
$ws.Range("J39").Value = 808.7778
$ws.Range("K39").Value = 255
$ws.Range("L39").Value = 2426.3334
$ws.Range("M39").Value = 41
$ws.Range("N39").Value = -3018.3334
$ws.Range("H41").Value = 66666890
$ws.Range("I41").Value = 111111256
$ws.Range("K41").Value = 111111256
$ws.Range("M41").Value = -111110816
$ws.Range("H76").Value = 3718.5
$ws.Range("J76").Value = 3892.7144
$ws.Range("L76").Value = 3892.7144
$ws.Range("N76").Value = -4522.7144
$ws.Range("H79").Value = 3718.5
$ws.Range("J79").Value = 3892.7144
$ws.Range("L79").Value = 3892.7144
$ws.Range("N79").Value = -6076.7144
$ws.Range("H113").Value = 1886.9131
$ws.Range("I113").Value = 1700
$ws.Range("J113").Value = 1938.8334
$ws.Range("K113").Value = 1700
$ws.Range("L113").Value = 1938.8334
$ws.Range("M113").Value = 1554
$ws.Range("N113").Value = -8446.8334
$ws.Range("H138").Value = 3096.3408
$ws.Range("I138").Value = 1315.25
$ws.Range("J138").Value = 4114.107
$ws.Range("K138").Value = 3945.75
$ws.Range("L138").Value = 12342.321
$ws.Range("M138").Value = 1194.25
$ws.Range("N138").Value = -22622.321

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3097.9333
$ws.Range("I45").Value = 3500
$ws.Range("J45").Value = 2829.889
$ws.Range("K45").Value = 3500
$ws.Range("L45").Value = 2829.889
$ws.Range("M45").Value = -3123
$ws.Range("N45").Value = -3583.889
$ws.Range("H53").Value = 333337660
$ws.Range("I53").Value = 333337660
$ws.Range("K53").Value = 333337660
$ws.Range("M53").Value = -333336978
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H61").Value = 2782.5
$ws.Range("I61").Value = 2707.6667
$ws.Range("J61").Value = 3007
$ws.Range("K61").Value = 2707.6667
$ws.Range("L61").Value = 3007
$ws.Range("M61").Value = -2495.6667
$ws.Range("N61").Value = -3431
$ws.Range("H136").Value = 2782.5
$ws.Range("I136").Value = 2707.6667
$ws.Range("J136").Value = 3007
$ws.Range("K136").Value = 8123.000100000001
$ws.Range("L136").Value = 9021
$ws.Range("M136").Value = -5573.000100000001
$ws.Range("N136").Value = -14121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 2000
$ws.Range("K29").Value = 2000
$ws.Range("M29").Value = -1711
$ws.Range("H56").Value = 11750
$ws.Range("I56").Value = 12000
$ws.Range("J56").Value = 11500
$ws.Range("K56").Value = 12000
$ws.Range("L56").Value = 11500
$ws.Range("M56").Value = -11261
$ws.Range("N56").Value = -12978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4157.375
$ws.Range("I16").Value = 1860
$ws.Range("K16").Value = 1860
$ws.Range("M16").Value = -1573
$ws.Range("H86").Value = 23929.182
$ws.Range("J86").Value = 29152.5
$ws.Range("L86").Value = 29152.5
$ws.Range("N86").Value = -31398.5
$ws.Range("H89").Value = 23929.182
$ws.Range("J89").Value = 29152.5
$ws.Range("L89").Value = 145762.5
$ws.Range("N89").Value = -156994.5
$ws.Range("H94").Value = 3385.739
$ws.Range("J94").Value = 4592.7856
$ws.Range("L94").Value = 4592.7856
$ws.Range("N94").Value = -5494.7856
$ws.Range("H99").Value = 2604
$ws.Range("I99").Value = 2044.4445
$ws.Range("J99").Value = 2918.75
$ws.Range("K99").Value = 2044.4445
$ws.Range("L99").Value = 2918.75
$ws.Range("M99").Value = -546.4445000000001
$ws.Range("N99").Value = -5914.75
$ws.Range("H107").Value = 1153.25
$ws.Range("I107").Value = 1100
$ws.Range("J107").Value = 1206.5
$ws.Range("K107").Value = 1100
$ws.Range("L107").Value = 1206.5
$ws.Range("M107").Value = 820
$ws.Range("N107").Value = -5046.5
$ws.Range("H113").Value = 4157.375
$ws.Range("I113").Value = 1860
$ws.Range("K113").Value = 1860
$ws.Range("M113").Value = 310
$ws.Range("H122").Value = 952.75
$ws.Range("I122").Value = 970.6667
$ws.Range("J122").Value = 899
$ws.Range("K122").Value = 2912.0001
$ws.Range("L122").Value = 2697
$ws.Range("M122").Value = -462.0001000000002
$ws.Range("N122").Value = -7597
$ws.Range("H126").Value = 2604
$ws.Range("I126").Value = 2044.4445
$ws.Range("J126").Value = 2918.75
$ws.Range("K126").Value = 6133.333500000001
$ws.Range("L126").Value = 8756.25
$ws.Range("M126").Value = -3663.333500000001
$ws.Range("N126").Value = -13696.25
$ws.Range("H134").Value = 2911.976
$ws.Range("I134").Value = 1672.8846
$ws.Range("J134").Value = 4925.5
$ws.Range("K134").Value = 5018.6538
$ws.Range("L134").Value = 14776.5
$ws.Range("M134").Value = -2483.6538
$ws.Range("N134").Value = -19846.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1463.8
$ws.Range("I4").Value = 439.66666
$ws.Range("K4").Value = 1318.99998
$ws.Range("M4").Value = -1206.99998
$ws.Range("H5").Value = 531.6896400000001
$ws.Range("J5").Value = 581.2727
$ws.Range("L5").Value = 1743.8181
$ws.Range("N5").Value = -1967.8181
$ws.Range("H107").Value = 38462096
$ws.Range("I107").Value = 264.5263
$ws.Range("J107").Value = 142858500
$ws.Range("K107").Value = 793.5789
$ws.Range("L107").Value = 428575500
$ws.Range("M107").Value = 1126.4211
$ws.Range("N107").Value = -428579340
$ws.Range("H132").Value = 594904.5
$ws.Range("I132").Value = 438.57144
$ws.Range("J132").Value = 1011030.7
$ws.Range("K132").Value = 3947.14296
$ws.Range("L132").Value = 9099276.299999999
$ws.Range("M132").Value = -1417.14296
$ws.Range("N132").Value = -9104336.299999999
$ws.Range("H135").Value = 531.6896400000001
$ws.Range("J135").Value = 581.2727
$ws.Range("L135").Value = 5231.454299999999
$ws.Range("N135").Value = -10301.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2731.8823
$ws.Range("I126").Value = 1951.7142
$ws.Range("J126").Value = 3278
$ws.Range("K126").Value = 5855.142599999999
$ws.Range("L126").Value = 9834
$ws.Range("M126").Value = -3385.142599999999
$ws.Range("N126").Value = -14774
$ws.Range("H132").Value = 5052.1113
$ws.Range("I132").Value = 5135.933
$ws.Range("J132").Value = 4633
$ws.Range("K132").Value = 15407.799
$ws.Range("L132").Value = 13899
$ws.Range("M132").Value = -12877.799
$ws.Range("N132").Value = -18959

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4114.5293
$ws.Range("I132").Value = 4414.8
$ws.Range("J132").Value = 3685.5715
$ws.Range("K132").Value = 13244.4
$ws.Range("L132").Value = 11056.7145
$ws.Range("M132").Value = -10714.4
$ws.Range("N132").Value = -16116.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2893.2
$ws.Range("I122").Value = 2280.3
$ws.Range("J122").Value = 3301.8
$ws.Range("K122").Value = 6840.900000000001
$ws.Range("L122").Value = 9905.400000000001
$ws.Range("M122").Value = -4390.900000000001
$ws.Range("N122").Value = -14805.4
